$d = $word.ActiveDocument

$replacements = @(
    @("25÷8=", "23÷8="),
    @("37÷9=", "71÷2="),
    @("17÷7=", "25÷5="),
    @("54÷5=", "33÷8="),
    @("90÷6=", "54÷7="),
    @("65÷2=", "60÷4="),
    @("51÷5=", "63÷3="),
    @("52÷9=", "41÷9="),
    @("23÷5=", "89÷5="),
    @("51÷8=", "43÷3="),
    @("24÷8=", "50÷3="),
    @("94÷9=", "86÷4="),
    @("99÷6=", "55÷3="),
    @("33÷9=", "18÷8="),
    @("95÷9=", "22÷4="),
    @("45÷6=", "46÷3="),
    @("55÷8=", "93÷6="),
    @("96÷7=", "14÷6="),
    @("87÷7=", "24÷5="),
    @("75÷2=", "56÷6="),
    @("57÷9=", "38÷5="),
    @("30÷5=", "26÷9="),
    @("27÷3=", "99÷8="),
    @("44÷7=", "73÷7="),
    @("61÷9=", "44÷9=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
